$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 19 data (string cells set in shared-string insertion order)
$ws.Range("A19").Value = 242
$ws.Range("B19").Value = "Valid Anagram"
$ws.Range("E19").Value = "Contain char c and Remove char c"
$ws.Range("C19").Value = "Frequency Table/List/LinQ/HashMap"
$ws.Range("D19").Value = "Frequency table ++ --, 2 HashMap ++  and compare them"

# Add new row 20 data
$ws.Range("A20").Value = 3330
$ws.Range("B20").Value = "Find the original typed string"
$ws.Range("C20").Value = "Frequency Table/Skip/Compare prev"

# Update the view/selection to match new state
$ws.Range("D13").Select()
